$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Demand")

# Extend the time-series (column A: index 2..12, column B: demand value)
# row 3 already holds index 1; add rows 4..14 for indices 2..12.
for ($i = 4; $i -le 14; $i++) {
    $ws.Cells.Item($i, 1).Value = $i - 2
}

# Update the demand value for the existing row and propagate it down the
# newly added rows (same constant value for every time step).
$ws.Range("B3:B14").Value = 497416667

# Copy B3's formatting (style s="9") down into the new cells so the whole
# column keeps a consistent number format.
$ws.Range("B3").Copy()
$ws.Range("B4:B14").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Best-fit column B so it comfortably shows the larger numbers.
$ws.Columns.Item(2).ColumnWidth = 10.1666666666667

# Make "Demand" the active/visible sheet with F14 selected (matches the
# new selection left behind after editing the table).
$ws.Activate()
$ws.Range("F14").Select()
